$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Update HG23 (rows 22-26) and RS23 (row 27) data ---

# Row 22 (HG23-1 equivalent): E and G change, B-G get new font style
$ws.Range("E22").Value = 319.32
$ws.Range("G22").Value = 88.12

# Row 23: D and F change
$ws.Range("D23").Value = 311.66
$ws.Range("F23").Value = 116.73

# Row 24: no value changes, only style

# Row 25: D and E change
$ws.Range("D25").Value = 296.44
$ws.Range("E25").Value = 289

# Row 26: D and E change
$ws.Range("D26").Value = 244.33
$ws.Range("E26").Value = 244.85

# Row 27: C and G change
$ws.Range("C27").Value = 7
$ws.Range("G27").Value = 65.53

# Apply new font (Calibri 11, black) + new style to B22:G27
$rng = $ws.Range("B22:G27")
$rng.Font.Color = 0

# --- View changes ---
$ws.Activate()
$excel.ActiveWindow.Zoom = 85
$ws.Range("B11").Select()
